# Add 2022-Q3 data
# ------------------------------------------------------------------
# The workbook tracks quarterly fund-holdings snapshots, one sheet per
# quarter, plus a "总计" (totals) roll-up sheet. A new quarter, 2022-Q3,
# is being inserted as the most-recent snapshot:
#   - a new "2022-Q3" worksheet is inserted right after "总计" (pushing
#     2022-Q2 / 2022-Q1 / 2021-Q4 down by one tab each, unchanged)
#   - the "总计" sheet gets a new row for 2022-Q3 or its existing rows
#     shift down to make room
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" (totals) sheet: insert a row for the new 2022-Q3 quarter.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push rows 3 down to make room for the new "2022-Q1" row, and extend
# with a freshly-formatted row 5 for "2021-Q4" (copies formatting from
# the existing data row).
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A3:D3").Copy($total.Range("A5:D5"))

# Row 2: 2022-Q3 (new)
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.01

# Row 3: 2022-Q2 (was row 2)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.01

# Row 4: 2022-Q1 (was row 3)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.04

# Row 5: 2021-Q4 (was row 4)
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 4
$total.Range("D5").Value = 0.06

# ------------------------------------------------------------------
# 2) Insert the new "2022-Q3" sheet right after "总计". Copying the
#    existing "2022-Q2" sheet brings along matching column widths /
#    header styling, and automatically pushes 2022-Q2 / 2022-Q1 /
#    2021-Q4 one tab to the right, unchanged.
# ------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $total)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Extend from 2 data rows to 4 data rows (copy row3's formatting down).
$q3.Range("A3:H3").Copy($q3.Range("A4:H4"))
$q3.Range("A3:H3").Copy($q3.Range("A5:H5"))

# Helper pattern for text-typed cells (fund code / name / percentages
# are stored as text, even when they look numeric, e.g. "013166"):
# set NumberFormat to Text first so values like "013166" or "0.38"
# aren't coerced to numbers, then reset the style back to Normal so no
# visible formatting change remains on the cell.

# Row 2: 013166 东兴宸祥量化混合A
$rng = $q3.Range("B2:G2")
$rng.NumberFormat = "@"
$q3.Range("B2").Value = "013166"
$q3.Range("C2").Value = "东兴宸祥量化混合A"
$q3.Range("D2").Value = "0.38"
$q3.Range("E2").Value = "93.87"
$q3.Range("F2").Value = "1.28"
$q3.Range("G2").Value = "0.0049"
$rng.Style = "Normal"
$q3.Range("H2").Value = 2

# Row 3: 009327 东兴兴晟混合A
$rng = $q3.Range("B3:G3")
$rng.NumberFormat = "@"
$q3.Range("B3").Value = "009327"
$q3.Range("C3").Value = "东兴兴晟混合A"
$q3.Range("D3").Value = "0.38"
$q3.Range("E3").Value = "79.70"
$q3.Range("F3").Value = "1.16"
$q3.Range("G3").Value = "0.0044"
$rng.Style = "Normal"
$q3.Range("H3").Value = 2

# Row 4: 013167 东兴宸祥量化混合C
$q3.Range("A4").Value = 2
$rng = $q3.Range("B4:G4")
$rng.NumberFormat = "@"
$q3.Range("B4").Value = "013167"
$q3.Range("C4").Value = "东兴宸祥量化混合C"
$q3.Range("D4").Value = "0.08"
$q3.Range("E4").Value = "93.87"
$q3.Range("F4").Value = "1.28"
$q3.Range("G4").Value = "0.0010"
$rng.Style = "Normal"
$q3.Range("H4").Value = 2

# Row 5: 009328 东兴兴晟混合C
$q3.Range("A5").Value = 3
$rng = $q3.Range("B5:G5")
$rng.NumberFormat = "@"
$q3.Range("B5").Value = "009328"
$q3.Range("C5").Value = "东兴兴晟混合C"
$q3.Range("D5").Value = "0.07"
$q3.Range("E5").Value = "79.70"
$q3.Range("F5").Value = "1.16"
$q3.Range("G5").Value = "0.0008"
$rng.Style = "Normal"
$q3.Range("H5").Value = 2
